$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("2_word_algorithm")

# Fill in the "options" indicator column (E) for rows 4-8 with "yes",
# matching the style already used in row 10 (s=4).
$ws.Range("E4:E8").Value = "yes"
$ws.Range("E10").Copy()
$ws.Range("E4:E8").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Move the cell selection/cursor to G4 (was L16).
$ws.Range("G4").Select()
